$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F13").Value = 1
$ws.Range("F18").Value = 0
$ws.Range("F19").Value = -1
$ws.Range("F22").Value = -1
$ws.Range("F27").Value = -1
$ws.Range("F29").Value = 2
$ws.Range("F38").Value = -1
$ws.Range("F49").Value = 1
$ws.Range("F51").Value = 0
$ws.Range("F55").Value = 2
$ws.Range("F56").Value = 1
$ws.Range("F63").Value = -6
$ws.Range("F65").Value = -3
$ws.Range("F69").Value = 13
